# Fix source properties generator script typo: add missing "comments" column
# to both data dictionary tables (one per worksheet).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $lo = $ws.ListObjects.Item(1)

    # Add a new trailing column to the table and give it a header
    $col = $lo.ListColumns.Add()
    $headerCell = $lo.HeaderRowRange.Columns.Item($lo.ListColumns.Count)
    $headerCell.Value = "comments"

    # Match the column width Excel would have used for the new column
    $newColIndex = $lo.Range.Column + $lo.Range.Columns.Count - 1
    $ws.Columns.Item($newColIndex).ColumnWidth = 11.498697916666666
}

# sheet2 ("data_src_a_table_b") also has a merged banner cell (H2:J2) above
# the table that needs to grow to keep spanning the header row
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("H2:K2").Merge()
